# Update the cryptocurrency price list with refreshed values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.957.84'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.36%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.863.13'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.72%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.26%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '336.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.52%  '

$ws.Range("E6").Value = '  -0.25%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4704'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.17%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3886'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.70%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.84'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.65%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07975'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.76%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9773'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.24%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.45'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.53%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.857.27'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.50%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.928'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.78%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.204'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.88%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.62'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.43%  '

$ws.Range("E17").Value = '  -0.30%  '

$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001037'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.09%  '

$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06626'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.46%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.44'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.60%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.003'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.25%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '27.940.89'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.32%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.389'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.28%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.90'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.03%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.293'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.77%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.066.69'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '158.94'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.05%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.55'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.04%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.094'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.23%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.446'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.47%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '119.31'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.57%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09477'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.67%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9559'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.41%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.579'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.19%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.305'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.06%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.344'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.24%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06086'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.59%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02233'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.03%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.288'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.36%  '

$ws.Range("E40").Value = '  -1.86%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.002'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.13%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5892'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.37%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1862'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.31%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.17'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.98%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.300'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.92%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.11'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.97%  '

$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5518'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.39%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.948'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.62%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06869'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.43%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '111.55'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.99%  '

$ws.Range("E51").Value = '  -32.80%  '
